# Agregado de modulo ABMUsuarios
# Inserts a new "FechaAlta" column between the existing "Mail" (F) and
# "Estado" (G) columns, pushing the old "Estado" column from G to H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the old column G ("Estado") content into the new column H ---
$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("H2").Value2 = $ws.Range("G2").Value2

# The "Dejar siempre en '1'" note currently lives on G1; it needs to move to
# H1 along with the values above. G1 already has a comment attached, so grab
# its text before it gets repurposed for the new column below, and create
# the note on H1.
$oldG1CommentText = $ws.Range("G1").Comment.Text()
$ws.Range("H1").AddComment($oldG1CommentText)

# --- 2. Turn column G into the new "FechaAlta" column ---
$ws.Range("G1").Value2 = "FechaAlta"
$null = $ws.Range("G1").Comment.Text("Fecha de alta en sistema con el siguiente formato:`nDD/MM/AAAA")

# Sample row: store the date as literal text (left aligned) rather than a
# real date serial number.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").HorizontalAlignment = -4131
$ws.Range("G2").Value2 = "1/4/2022"

# Leave the selection on the default top-left cell.
$null = $ws.Range("A1").Select()
